$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns at F:G, pushing the old "Topic" column (F) to H
$ws.Columns("F:G").Insert()

# New headers
$ws.Range("F1").Value = "Width"
$ws.Range("G1").Value = "Height"

# Fill in Width/Height values (12 x 12) for each data row
for ($r = 2; $r -le 12; $r++) {
    $ws.Cells.Item($r, 6).Value = 12
    $ws.Cells.Item($r, 7).Value = 12
}

# Row heights / default row height adjustments
$ws.Rows("1").RowHeight = 13.8

# Column width tweaks
$ws.Columns("A").ColumnWidth = 14.28
$ws.Columns("E").ColumnWidth = 12.1

# Update the active selection
$ws.Range("G2").Select()
